$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 736, shifting existing rows 736:779 down to 741:784.
$ws.Range("A736:T740").EntireRow.Insert()

# New data rows (dates are Excel serials in the 1900 date system).
$newRows = @(
    @(3,"Femacal de La Calera","Coquimbo",44714,5,"Fruta",100109,"Uva",100109001,"Uva","Autumn Royal","Primera",50,9000,9000,9000,"`$/caja 15 kilos","Provincia de San Felipe de Aconcagua",600,15),
    @(3,"Femacal de La Calera","Coquimbo",44714,5,"Fruta",100109,"Uva",100109001,"Uva","Crimpson Seedless","Primera",60,9000,9000,9000,"`$/caja 15 kilos","Provincia de San Felipe de Aconcagua",600,15),
    @(3,"Femacal de La Calera","Coquimbo",44714,5,"Fruta",100109,"Uva",100109001,"Uva","Red Globe","Primera",85,9000,9000,9000,"`$/caja 15 kilos","Provincia de San Felipe de Aconcagua",600,15),
    @(3,"Femacal de La Calera","Coquimbo",44714,5,"Fruta",100109,"Uva",100109001,"Uva","Rosada pastilla","Primera",70,16000,16000,16000,"`$/caja 15 kilos","Provincia de San Felipe de Aconcagua",1067,15),
    @(3,"Femacal de La Calera","Coquimbo",44714,5,"Fruta",100109,"Uva",100109001,"Uva","Thompson seedless","Primera",60,16000,16000,16000,"`$/caja 15 kilos","Provincia de San Felipe de Aconcagua",1067,15)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 736 + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    # Match the date-formatted style used by the existing "Fecha" column (D).
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
